$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")

# ---------------------------------------------------------------------------
# The "Settings" sheet had a copy/paste bug: row 16 only held the
# "ConnectionString" label in column A (column B was empty), while row 17
# duplicated that same label together with the actual connection-string
# value that belonged in row 16. Deleting the stray row 16 removes the
# duplicate and shifts every following row up by one (carrying their
# formatting with them, which is exactly what happened in the real edit).
# ---------------------------------------------------------------------------
$ws.Rows("16:16").Delete()

# Implement config base retry + global exception handler: two new rows are
# appended below the existing (now shifted-up) settings block.
$ws.Range("A21").Value = "GlobalExceptionMaxRetryCount"
$ws.Range("B21").Value = 3

$ws.Range("A22").Value = "GlobalExceptionRetryDelayInSec"
$ws.Range("B22").Value = 5

# Restore the selection/active-cell state captured in the workbook: the
# Constants sheet remembers A15 as its last selection, while the Settings
# sheet stays the active tab with C21 selected.
$wsConstants.Activate()
$wsConstants.Range("A15").Select()

$ws.Activate()
$ws.Range("C21").Select()
